$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New duty-roster values for column B (rows 2..32), matching the target state.
# Empty string clears the cell (same on-disk result as a blank <c s="2"/> cell).
$values = @{
    2  = ""
    3  = "Hansen Jakob U"
    4  = ""
    5  = "Nicholas Tristan Aryasatyo"
    6  = ""
    7  = "渡部魁"
    8  = "氏家琉貴"
    9  = "Yunjae"
    10 = "遠藤隼人"
    11 = "富澤天音"
    12 = ""
    13 = "川田涼介"
    14 = "志塚惇希"
    15 = "川田涼介"
    16 = "豊島亮"
    17 = "兒島大志郎"
    18 = "日高泰聖"
    19 = "白岩詩佑介"
    20 = "Cox Matthew Jonah"
    21 = "Hansen Jakob U"
    22 = "石井海成"
    23 = "Nicholas Tristan Aryasatyo"
    24 = "小溝賢"
    25 = "小野文哉"
    26 = "渡部魁"
    27 = "崎谷航平"
    28 = "三神佳誠"
    29 = "氏家琉貴"
    30 = "羽賀尚生"
    31 = "島田実"
    32 = "足立耕平"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}

# Update the active cell selection to match the saved view state.
[void]$ws.Range("D14").Select()
